$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "bugs" variant tags text in E11 and replace with numeric 0,
# which also drops the now-unused shared string from the workbook.
$ws.Range("E11").Value = 0

# Update the selection shown when the sheet is reopened.
$ws.Range("D16").Select()
